$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Change: "A auditoria agendada para 20/05/2015 " -> "A auditoria agendada "
# The run that holds this text sits immediately next to another run
# ("para o ") that has identical run formatting. Editing text in-place
# normally causes the engine to coalesce the two adjacent, identically
# formatted runs into a single run - which would alter document
# structure beyond what the source diff describes. To avoid that, we
# briefly give the neighboring run a distinguishing direct-formatting
# value (Bold) so it will not be merged while we edit the first run,
# then clear that temporary formatting back to "undefined" (which
# removes the direct formatting entirely, restoring the original
# appearance) once the text edit is complete.
# --------------------------------------------------------------------

$wdUndefined = 9999999

# Step 1: temporarily mark the neighboring run ("para o ") so it keeps
# its own identity while we edit the preceding run.
$guard = $d.Content
$guardFind = $guard.Find
$guardFind.ClearFormatting()
$guardFind.Execute("para o ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($guardFind.Found) {
    $guard.Bold = 1
}

# Step 2: replace the text of the target run.
$target = $d.Content
$targetFind = $target.Find
$targetFind.ClearFormatting()
$targetFind.Execute("A auditoria agendada para 20/05/2015 ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($targetFind.Found) {
    $target.Text = "A auditoria agendada "
}

# Step 3: remove the temporary guard formatting so the neighboring run
# returns to its original, unformatted state.
$unguard = $d.Content
$unguardFind = $unguard.Find
$unguardFind.ClearFormatting()
$unguardFind.Execute("para o ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($unguardFind.Found) {
    $unguard.Bold = $wdUndefined
}
